# Actualización desde MV -datos-
# Add the three new daily auction rows (05-10-2021, 06-10-2021, 07-10-2021)
# below the existing data table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("05-10-2021", "06-10-2021", "07-10-2021")
$values = @(
    @(40, 191, 40, 813),
    @(40, 162, 40, 816),
    @(40, 163, 40, 814)
)

$startRow = 182
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i

    # Force column A to be stored as literal text (matching the existing
    # "dd-mm-yyyy" string cells above) instead of being auto-parsed as a date,
    # then drop back to the sheet's default (unstyled) formatting.
    $ws.Range("A$row").NumberFormat = "@"
    $ws.Range("A$row").Value = $dates[$i]
    $ws.Range("A$row").ClearFormats()

    $ws.Range("B$row").Value = $values[$i][0]
    $ws.Range("C$row").Value = $values[$i][1]
    $ws.Range("D$row").Value = $values[$i][2]
    $ws.Range("E$row").Value = $values[$i][3]
}
